# Append the new challenge row (row 3 in the sheet's existing 0-based row
# numbering, i.e. the 4th physical row after header row "1" and data rows
# "0"/"2") with the same column layout as the existing rows:
#   A: gameID        B: challenger   C: rating   D: wager
#   E: link          F: escrowID     G: accepted?   H: accepter
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ABxocmbB"
$ws.Range("B3").Value = "trashboatsr"
$ws.Range("C3").Value = 1818
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = "https://lichess.org/ABxocmbB"
$ws.Range("F3").Value = 3024
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = "blank"
